$d = $word.ActiveDocument

# --- 1. Split the trailing period off of the "mine the block." run into its own run ---
# Locate the paragraph ending in "...mine the block."
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*mine the block.*") {
        $targetPara = $p
    }
}

$rngAll = $targetPara.Range
$periodPos = $rngAll.End - 1
$period = $d.Range($periodPos - 1, $periodPos)
$fmtCopy = $period.FormattedText

# Insert a formatted copy of the period just before the original period
# (this creates a new run carrying the same rPr as the source run).
$insertPoint = $d.Range($periodPos - 1, $periodPos - 1)
$insertPoint.FormattedText = $fmtCopy

# Delete the now-duplicated trailing period, leaving the newly inserted one
# as its own separate run.
$pNow = $targetPara
$rngNow = $pNow.Range
$lastPeriod = $d.Range($rngNow.End - 2, $rngNow.End - 1)
$lastPeriod.Delete()

# --- 2. Remove the trailing "still need to pretty this up" scratch paragraphs ---
# Keep the paragraph that holds the _GoBack bookmark, but drop the empty
# paragraph after it and the "STILL NEED TO PRETTY THIS UP..." paragraph
# that follows.
$scratchPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*STILL NEED TO PRETTY THIS UP*") {
        $scratchPara = $p
    }
}
$scratchIndex = $scratchPara.Index
$emptyPara = $d.Paragraphs.Item($scratchIndex - 1)

# Delete the empty paragraph first (it is not the last paragraph in the
# document, so its paragraph mark can be removed outright).
$emptyPara.Range.Delete()

# The scratch paragraph is now the last paragraph in the document; its
# paragraph mark can't be deleted directly, but its text content can -
# leaving the previous (bookmark) paragraph as the document's last
# paragraph once its own (now-empty) mark collapses away.
$scratchPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratchPara2.Range.Delete()
